$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8:I8").Value = "UC09a_TS_001"
$ws.Range("B10").Value = "TC_Ser_001"
$ws.Range("C10").Value = "TC_Ser_002"
$ws.Range("D10").Value = "TC_Ser_003"
$ws.Range("E10").Value = "TC_Ser_004"
$ws.Range("F10").Value = "TC_Ser_005"
$ws.Range("G10").Value = "TC_Ser_006"
$ws.Range("H10").Value = "TC_Ser_007"
$ws.Range("I10").Value = "TC_Ser_008"
